$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row 67 with the new "district heating demand" input field
$ws.Range("A67").Value = "district heating demand "
$ws.Range("B67").Value = "input"
$ws.Range("C67").Value = "dh_demand"
$ws.Range("D67").Value = 100
$ws.Range("E67").Value = 0
$ws.Range("F67").Value = "MWh"
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 9.9999999999999998E+101

# View / selection updates to match the target sheet view
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("F70").Select()
